$d = $word.ActiveDocument
$wns = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

# ---------------------------------------------------------------------------
# Starting layout:
#   Paragraph 1: "Prueba fca"           (spell-checked run pair, unchanged)
#   Paragraph 2: bookmark "_GoBack"     (unchanged)
#
# Target layout:
#   Paragraph 1: "Prueba fca"           (unchanged)
#   Paragraph 2: <w:p/>                 (new, empty)
#   Paragraph 3: "Prueba efu"           (new, spell-checked run pair)
#   Paragraph 4: bookmark "_GoBack"     (unchanged)
#   Paragraph 5: <w:p/>                 (new, empty)
# ---------------------------------------------------------------------------

# 1) Insert a new (currently empty) paragraph right after paragraph 1.
$p1 = $d.Paragraphs.Item(1)
$p1.Range.InsertParagraphAfter()

# Normalize it to a bare, truly empty paragraph element (<w:p/>).
$p2 = $d.Paragraphs.Item(2)
$p2full = $d.Range($p2.Range.Start, $p2.Range.End)
$null = $p2full.InsertXML("<w:p $wns/>")

# 2) Insert another new paragraph right after that empty one - this will
#    become the "Prueba efu" paragraph.
$p2 = $d.Paragraphs.Item(2)
$p2.Range.InsertParagraphAfter()

$p3 = $d.Paragraphs.Item(3)
$p3full = $d.Range($p3.Range.Start, $p3.Range.End)
$p3xml = "<w:p $wns>" +
         "<w:r><w:t xml:space='preserve'>Prueba </w:t></w:r>" +
         "<w:proofErr w:type='spellStart'/>" +
         "<w:r><w:t>efu</w:t></w:r>" +
         "<w:proofErr w:type='spellEnd'/>" +
         "</w:p>"
$null = $p3full.InsertXML($p3xml)

# 3) Insert a trailing empty paragraph right after the bookmark paragraph
#    (paragraph 4), before the final section break.
$p4 = $d.Paragraphs.Item(4)
$p4.Range.InsertParagraphAfter()

$p5 = $d.Paragraphs.Item(5)
$p5full = $d.Range($p5.Range.Start, $p5.Range.End)
$null = $p5full.InsertXML("<w:p $wns/>")
